$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formats first (reusing existing style indices via copy/paste of formats) ---

# Row 5 header formatting: reuse the yellow header style (style index 1) from row 2.
$ws.Range("A2:C2").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)

# Row 6 data formatting: reuse the bordered body style (style index 2) from row 3 (A3:B3, not C3 which is date-formatted).
$ws.Range("A3:B3").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("F6").PasteSpecial(-4122)

# E6 gets the bordered style too, then wrap text is turned on (creates the new wrap-text style).
$ws.Range("A3").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").WrapText = $true

# --- Now fill in values (order matters for shared-string table layout) ---
$ws.Range("A5").Value = "Opportunity_Name"
$ws.Range("B5").Value = "Sales_Stage"
$ws.Range("C5").Value = "Amount"
$ws.Range("D5").Value = "Assigned_To"
$ws.Range("E5").Value = "Description"

$ws.Range("A6").Value = "Ross Taylor"
$ws.Range("B6").Value = "Value Proposition"
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = "Support Group"
$ws.Range("E6").Value = "Ross Taylor is kiwi business man. `nHe is investing his money in the `nnew business."

$ws.Range("F5").Value = "Organisation_Name"
$ws.Range("F6").Value = "Kiwi organisations"

# --- Row height / column widths ---
$ws.Rows.Item(6).RowHeight = 60

$ws.Columns.Item(1).ColumnWidth = 17.59
$ws.Columns.Item(2).ColumnWidth = 16.09
$ws.Columns.Item(3).ColumnWidth = 7.25
$ws.Columns.Item(4).ColumnWidth = 13.09
$ws.Columns.Item(5).ColumnWidth = 28.76
$ws.Columns.Item(6).ColumnWidth = 17.92

# --- Selection matches the last-edited cell ---
[void]$ws.Range("F6").Select()
